$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('LP1912')
$ws.Cells.Item(2, 1).Value = 'Última actualización: 21:58:04'
$ws.Cells.Item(3, 1).Value = 'Total filas: 572'
$ws.Cells.Item(48, 1).Value = '06:02:16'
$ws.Cells.Item(48, 3).Value = '23_HERNANDEZ'
$ws.Cells.Item(48, 4).Value = 63
$ws.Cells.Item(49, 1).Value = '05:18:23'
$ws.Cells.Item(49, 3).Value = '15_ABASTO'
$ws.Cells.Item(49, 4).Value = 107
$ws.Cells.Item(84, 1).Value = '06:37:24'
$ws.Cells.Item(84, 3).Value = '215B_EL PATO'
$ws.Cells.Item(84, 4).Value = 106
$ws.Cells.Item(85, 1).Value = '07:44:08'
$ws.Cells.Item(85, 3).Value = '16_P MOR-SANTA ANA'
$ws.Cells.Item(85, 4).Value = 39
$ws.Cells.Item(91, 1).Value = '07:14:27'
$ws.Cells.Item(91, 3).Value = '17_ROMERO'
$ws.Cells.Item(91, 4).Value = 100
$ws.Cells.Item(92, 1).Value = '08:47:19'
$ws.Cells.Item(92, 3).Value = '10_OLMOS'
$ws.Cells.Item(92, 4).Value = 7
$ws.Cells.Item(111, 1).Value = '08:47:19'
$ws.Cells.Item(111, 3).Value = '16_SANTA ANA'
$ws.Cells.Item(111, 4).Value = 36
$ws.Cells.Item(113, 1).Value = '07:57:27'
$ws.Cells.Item(113, 3).Value = '11_ETCHEVERRY'
$ws.Cells.Item(113, 4).Value = 86
$ws.Cells.Item(120, 3).Value = '23_HERNANDEZ'
$ws.Cells.Item(121, 3).Value = '16_SANTA ANA'
$ws.Cells.Item(183, 1).Value = '11:43:19'
$ws.Cells.Item(183, 3).Value = '16_SANTA ANA'
$ws.Cells.Item(183, 4).Value = 0
$ws.Cells.Item(184, 1).Value = '10:50:37'
$ws.Cells.Item(184, 3).Value = '10_OLMOS'
$ws.Cells.Item(184, 4).Value = 53
$ws.Cells.Item(185, 3).Value = '17_ROMERO'
$ws.Cells.Item(272, 1).Value = '13:24:27'
$ws.Cells.Item(272, 3).Value = '11_ETCHEVERRY'
$ws.Cells.Item(272, 4).Value = 41
$ws.Cells.Item(273, 1).Value = '12:44:05'
$ws.Cells.Item(273, 3).Value = '23_HERNANDEZ'
$ws.Cells.Item(273, 4).Value = 81
$ws.Cells.Item(280, 1).Value = '13:51:56'
$ws.Cells.Item(280, 3).Value = '26_HERNANDEZ'
$ws.Cells.Item(280, 4).Value = 29
$ws.Cells.Item(281, 1).Value = '12:24:14'
$ws.Cells.Item(281, 3).Value = '215C_EL PATO'
$ws.Cells.Item(281, 4).Value = 116
$ws.Cells.Item(316, 1).Value = '13:51:56'
$ws.Cells.Item(316, 3).Value = '215A_EL PATO'
$ws.Cells.Item(316, 4).Value = 107
$ws.Cells.Item(317, 1).Value = '14:17:27'
$ws.Cells.Item(317, 3).Value = '23_HERNANDEZ'
$ws.Cells.Item(317, 4).Value = 81
$ws.Cells.Item(326, 3).Value = '11_ETCHEVERRY'
$ws.Cells.Item(327, 3).Value = '17_ROMERO'
$ws.Cells.Item(335, 1).Value = '15:21:47'
$ws.Cells.Item(335, 3).Value = '16_SANTA ANA'
$ws.Cells.Item(335, 4).Value = 41
$ws.Cells.Item(336, 1).Value = '14:56:04'
$ws.Cells.Item(336, 3).Value = '10_OLMOS'
$ws.Cells.Item(336, 4).Value = 66
$ws.Cells.Item(362, 3).Value = '225_GOMEZ'
$ws.Cells.Item(363, 1).Value = '14:56:04'
$ws.Cells.Item(363, 3).Value = '16_P MOR-SANTA ANA'
$ws.Cells.Item(363, 4).Value = 107
$ws.Cells.Item(364, 1).Value = '16:14:21'
$ws.Cells.Item(364, 3).Value = '23_HERNANDEZ'
$ws.Cells.Item(364, 4).Value = 29
$ws.Cells.Item(380, 1).Value = '16:39:47'
$ws.Cells.Item(380, 3).Value = '23_HERNANDEZ'
$ws.Cells.Item(380, 4).Value = 28
$ws.Cells.Item(381, 1).Value = '16:30:20'
$ws.Cells.Item(381, 3).Value = '27_EL RETIRO'
$ws.Cells.Item(381, 4).Value = 37
$ws.Cells.Item(467, 1).Value = '18:34:43'
$ws.Cells.Item(467, 3).Value = '15_ABASTO'
$ws.Cells.Item(467, 4).Value = 42
$ws.Cells.Item(468, 1).Value = '17:34:37'
$ws.Cells.Item(468, 3).Value = '17_ROMERO'
$ws.Cells.Item(468, 4).Value = 102
$ws.Cells.Item(469, 3).Value = '27_EL RETIRO'
$ws.Cells.Item(485, 1).Value = '19:13:07'
$ws.Cells.Item(485, 3).Value = '17_ROMERO'
$ws.Cells.Item(485, 4).Value = 26
$ws.Cells.Item(487, 1).Value = '19:38:38'
$ws.Cells.Item(487, 3).Value = '16_SANTA ANA'
$ws.Cells.Item(487, 4).Value = 1
$ws.Cells.Item(505, 1).Value = '18:34:43'
$ws.Cells.Item(505, 3).Value = '16_P MOR-167 Y 521'
$ws.Cells.Item(505, 4).Value = 96
$ws.Cells.Item(506, 1).Value = '19:38:38'
$ws.Cells.Item(506, 3).Value = '10_OLMOS'
$ws.Cells.Item(506, 4).Value = 32
$ws.Cells.Item(532, 3).Value = '16_SANTA ANA'
$ws.Cells.Item(533, 3).Value = '17_ROMERO'
$ws.Cells.Item(534, 3).Value = '23_HERNANDEZ'
$ws.Cells.Item(540, 1).Value = '19:13:07'
$ws.Cells.Item(540, 3).Value = '84_COLONIA URQUIZA-ESC 49'
$ws.Cells.Item(540, 4).Value = 111
$ws.Cells.Item(541, 1).Value = '20:33:25'
$ws.Cells.Item(541, 3).Value = '15_ABASTO'
$ws.Cells.Item(541, 4).Value = 31
$ws.Cells.Item(558, 1).Value = '21:58:04'
$ws.Cells.Item(558, 3).Value = '17_ROMERO'
$ws.Cells.Item(558, 4).Value = 10
$ws.Cells.Item(559, 2).Value = '22:08'
$ws.Cells.Item(559, 3).Value = '11_ETCHEVERRY'
$ws.Cells.Item(559, 4).Value = 95
$ws.Cells.Item(560, 1).Value = '20:33:25'
$ws.Cells.Item(560, 2).Value = '22:19'
$ws.Cells.Item(560, 4).Value = 106
$ws.Cells.Item(561, 2).Value = '22:22'
$ws.Cells.Item(561, 3).Value = '26_HERNANDEZ'
$ws.Cells.Item(561, 4).Value = 96
$ws.Cells.Item(562, 1).Value = '21:58:04'
$ws.Cells.Item(562, 2).Value = '22:23'
$ws.Cells.Item(562, 3).Value = '23_HERNANDEZ'
$ws.Cells.Item(562, 4).Value = 25
$ws.Cells.Item(563, 1).Value = '20:46:33'
$ws.Cells.Item(563, 2).Value = '22:27'
$ws.Cells.Item(563, 3).Value = '84_COLONIA URQUIZA-ESC 49'
$ws.Cells.Item(563, 4).Value = 101
$ws.Cells.Item(564, 1).Value = '20:33:25'
$ws.Cells.Item(564, 2).Value = '22:28'
$ws.Cells.Item(564, 3).Value = '84_COLONIA URQUIZA-ESC 49'
$ws.Cells.Item(564, 4).Value = 115
$ws.Cells.Item(565, 1).Value = '21:58:04'
$ws.Cells.Item(565, 2).Value = '22:28'
$ws.Cells.Item(565, 3).Value = '10_OLMOS'
$ws.Cells.Item(565, 4).Value = 30
$ws.Cells.Item(566, 2).Value = '22:29'
$ws.Cells.Item(566, 3).Value = '10_OLMOS'
$ws.Cells.Item(566, 4).Value = 96
$ws.Cells.Item(567, 1).Value = '20:46:33'
$ws.Cells.Item(567, 2).Value = '22:30'
$ws.Cells.Item(567, 3).Value = '10_OLMOS'
$ws.Cells.Item(567, 4).Value = 104
$ws.Cells.Item(567, 5).Value = 'LP1912'
$ws.Cells.Item(568, 1).Value = '20:46:33'
$ws.Cells.Item(568, 2).Value = '22:39'
$ws.Cells.Item(568, 3).Value = '215A_EL PATO'
$ws.Cells.Item(568, 4).Value = 113
$ws.Cells.Item(568, 5).Value = 'LP1912'
$ws.Cells.Item(569, 1).Value = '21:58:04'
$ws.Cells.Item(569, 2).Value = '22:42'
$ws.Cells.Item(569, 3).Value = '26_HERNANDEZ'
$ws.Cells.Item(569, 4).Value = 44
$ws.Cells.Item(569, 5).Value = 'LP1912'
$ws.Cells.Item(570, 1).Value = '21:58:04'
$ws.Cells.Item(570, 2).Value = '22:46'
$ws.Cells.Item(570, 3).Value = '16_SANTA ANA'
$ws.Cells.Item(570, 4).Value = 48
$ws.Cells.Item(570, 5).Value = 'LP1912'
$ws.Cells.Item(571, 1).Value = '20:53:41'
$ws.Cells.Item(571, 2).Value = '22:50'
$ws.Cells.Item(571, 3).Value = '14_ABASTO'
$ws.Cells.Item(571, 4).Value = 117
$ws.Cells.Item(571, 5).Value = 'LP1912'
$ws.Cells.Item(572, 1).Value = '21:58:04'
$ws.Cells.Item(572, 2).Value = '22:52'
$ws.Cells.Item(572, 3).Value = '14_ABASTO'
$ws.Cells.Item(572, 4).Value = 54
$ws.Cells.Item(572, 5).Value = 'LP1912'
$ws.Cells.Item(573, 1).Value = '21:58:04'
$ws.Cells.Item(573, 2).Value = '23:08'
$ws.Cells.Item(573, 3).Value = '14X44_ABASTO'
$ws.Cells.Item(573, 4).Value = 70
$ws.Cells.Item(573, 5).Value = 'LP1912'
$ws.Cells.Item(574, 1).Value = '21:58:04'
$ws.Cells.Item(574, 2).Value = '23:08'
$ws.Cells.Item(574, 3).Value = '17_ROMERO'
$ws.Cells.Item(574, 4).Value = 70
$ws.Cells.Item(574, 5).Value = 'LP1912'
$ws.Cells.Item(575, 1).Value = '21:58:04'
$ws.Cells.Item(575, 2).Value = '23:12'
$ws.Cells.Item(575, 3).Value = '23_HERNANDEZ'
$ws.Cells.Item(575, 4).Value = 74
$ws.Cells.Item(575, 5).Value = 'LP1912'
$ws.Cells.Item(576, 1).Value = '21:58:04'
$ws.Cells.Item(576, 2).Value = '23:12'
$ws.Cells.Item(576, 3).Value = '16_SANTA ANA'
$ws.Cells.Item(576, 4).Value = 74
$ws.Cells.Item(576, 5).Value = 'LP1912'
$ws.Cells.Item(577, 1).Value = '21:58:04'
$ws.Cells.Item(577, 2).Value = '23:51'
$ws.Cells.Item(577, 3).Value = '215_ALUAR'
$ws.Cells.Item(577, 4).Value = 113
$ws.Cells.Item(577, 5).Value = 'LP1912'

$ws = $wb.Worksheets.Item('LP1912-215')
$ws.Cells.Item(2, 1).Value = 'Última actualización: 21:58:04'
$ws.Cells.Item(3, 1).Value = 'Total filas: 52'
$ws.Cells.Item(57, 1).Value = '21:58:04'
$ws.Cells.Item(57, 2).Value = '23:51'
$ws.Cells.Item(57, 3).Value = '215_ALUAR'
$ws.Cells.Item(57, 4).Value = 113
$ws.Cells.Item(57, 5).Value = 'LP1912'

$ws = $wb.Worksheets.Item('6203-6173')
$ws.Cells.Item(2, 1).Value = 'Última actualización: 21:58:04'
$ws.Cells.Item(3, 1).Value = 'Total filas: 70'
$ws.Cells.Item(74, 1).Value = '21:58:04'
$ws.Cells.Item(74, 2).Value = '22:34'
$ws.Cells.Item(74, 3).Value = '215B_LP-P MOR-40 Y 115'
$ws.Cells.Item(74, 4).Value = 36
$ws.Cells.Item(74, 5).Value = 'L6173'
$ws.Cells.Item(75, 1).Value = '21:58:04'
$ws.Cells.Item(75, 2).Value = '23:08'
$ws.Cells.Item(75, 3).Value = '215A_LA PLATA'
$ws.Cells.Item(75, 4).Value = 70
$ws.Cells.Item(75, 5).Value = 'L6173'

